$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'79.364.59"
$ws.Range('E2').Value = '  +3.98%  '

$ws.Range('D3').Value = "'3.165.69"
$ws.Range('E3').Value = '  +2.52%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = "'205.86"
$ws.Range('E5').Value = '  +3.89%  '

$ws.Range('D6').Value = "'626.44"
$ws.Range('E6').Value = '  +1.44%  '

$ws.Range('E7').Value = '  +28.05%  '

$ws.Range('D9').Value = "'0.589"
$ws.Range('E9').Value = '  +6.86%  '

$ws.Range('D10').Value = "'3.165.18"
$ws.Range('E10').Value = '  +2.56%  '

$ws.Range('D11').Value = "'0.590"
$ws.Range('E11').Value = '  +34.11%  '

$ws.Range('D12').Value = "'0.0000254"
$ws.Range('E12').Value = '  +31.06%  '

$ws.Range('E13').Value = '  +2.04%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = "'3.749.37"
$ws.Range('E14').Value = '  +2.60%  '

$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').Value = "'5.27"
$ws.Range('E15').Value = '  +0.74%  '

$ws.Range('D16').Value = "'31.51"
$ws.Range('E16').Value = '  +7.74%  '

$ws.Range('D17').Value = "'79.515.27"
$ws.Range('E17').Value = '  +3.97%  '

$ws.Range('D18').Value = "'3.156.12"
$ws.Range('E18').Value = '  +2.58%  '

$ws.Range('D19').Value = "'14.31"
$ws.Range('E19').Value = '  +5.21%  '

$ws.Range('D20').Value = "'437.03"
$ws.Range('E20').Value = '  +14.53%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'9.12"
$ws.Range('E21').Value = '  +0.41%  '

$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').Value = "'2.93"
$ws.Range('E22').Value = '  +13.60%  '

$ws.Range('D23').Value = "'5.23"
$ws.Range('E23').Value = '  +17.96%  '

$ws.Range('D24').Value = "'6.78"
$ws.Range('E24').Value = '  +5.06%  '

$ws.Range('D25').Value = "'3.331.77"
$ws.Range('E25').Value = '  +2.87%  '

$ws.Range('D26').Value = "'75.98"
$ws.Range('E26').Value = '  +5.06%  '

$ws.Range('D27').Value = "'4.66"
$ws.Range('E27').Value = '  +4.94%  '

$ws.Range('D28').Value = "'10.84"
$ws.Range('E28').Value = '  +8.15%  '

$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -0.51%  '

$ws.Range('D30').Value = "'0.0000121"
$ws.Range('E30').Value = '  +11.62%  '

$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = "'1.00"
$ws.Range('E31').Value = '  +0.19%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'9.02"
$ws.Range('E32').Value = '  +8.39%  '

$ws.Range('D33').Value = "'549.36"
$ws.Range('E33').Value = '  +9.68%  '

$ws.Range('D34').Value = "'1.47"
$ws.Range('E34').Value = '  +3.46%  '

$ws.Range('D35').Value = "'1.99"
$ws.Range('E35').Value = '  +3.87%  '

$ws.Range('D36').Value = "'0.148"
$ws.Range('E36').Value = '  +20.89%  '

$ws.Range('D37').Value = "'23.06"
$ws.Range('E37').Value = '  +10.86%  '

$ws.Range('D38').Value = "'0.123"
$ws.Range('E38').Value = '  +19.59%  '

$ws.Range('E39').Value = '  -0.04%  '

$ws.Range('D40').Value = "'0.407"
$ws.Range('E40').Value = '  +7.39%  '

$ws.Range('E41').Value = '  +3.43%  '

$ws.Range('D42').Value = "'164.12"
$ws.Range('E42').Value = '  +1.45%  '

$ws.Range('D43').Value = "'5.63"
$ws.Range('E43').Value = '  +9.72%  '

$ws.Range('D45').Value = "'187.74"
$ws.Range('E45').Value = '  -4.02%  '

$ws.Range('E46').Value = '  +8.99%  '

$ws.Range('D47').Value = "'2.67"
$ws.Range('E47').Value = '  +10.11%  '

$ws.Range('D48').Value = "'0.779"
$ws.Range('E48').Value = '  -3.40%  '

$ws.Range('D49').Value = "'1.30"
$ws.Range('E49').Value = '  +4.31%  '

$ws.Range('D50').Value = "'43.27"
$ws.Range('E50').Value = '  +4.86%  '

$ws.Range('D51').Value = "'4.26"
$ws.Range('E51').Value = '  +9.09%  '
